$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 14707200
$ws.Range("I43").Value = 45456210
$ws.Range("J43").Value = 1152.2609
$ws.Range("K43").Value = 45456210
$ws.Range("L43").Value = 1152.2609
$ws.Range("M43").Value = -45456141
$ws.Range("N43").Value = -1290.2609
# Row 58
$ws.Range("H58").Value = 928.5238000000001
$ws.Range("I58").Value = 85.57143000000001
$ws.Range("J58").Value = 1350
$ws.Range("K58").Value = 256.71429
$ws.Range("L58").Value = 4050
$ws.Range("M58").Value = -106.71429
$ws.Range("N58").Value = -4350
# Row 70
$ws.Range("H70").Value = 2255.5
$ws.Range("I70").Value = 2625.8333
$ws.Range("J70").Value = 1700
$ws.Range("K70").Value = 7877.499899999999
$ws.Range("L70").Value = 5100
$ws.Range("M70").Value = -7607.499899999999
$ws.Range("N70").Value = -5640
# Row 73
$ws.Range("H73").Value = 2255.5
$ws.Range("I73").Value = 2625.8333
$ws.Range("J73").Value = 1700
$ws.Range("K73").Value = 7877.499899999999
$ws.Range("L73").Value = 5100
$ws.Range("M73").Value = -6941.499899999999
$ws.Range("N73").Value = -6972
# Row 80
$ws.Range("H80").Value = 634.9
$ws.Range("I80").Value = 591.6667
$ws.Range("J80").Value = 699.75
$ws.Range("K80").Value = 1775.0001
$ws.Range("L80").Value = 2099.25
$ws.Range("M80").Value = -777.0001
$ws.Range("N80").Value = -4095.25
# Row 83
$ws.Range("H83").Value = 634.9
$ws.Range("I83").Value = 591.6667
$ws.Range("J83").Value = 699.75
$ws.Range("K83").Value = 5325.0003
$ws.Range("L83").Value = 6297.75
$ws.Range("M83").Value = -333.0002999999997
$ws.Range("N83").Value = -16281.75
# Row 103
$ws.Range("H103").Value = 1500
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Range("H31").Value = 2275
$ws.Range("I31").Value = 2275
$ws.Range("K31").Value = 2275
$ws.Range("M31").Value = -1981
# Row 32
$ws.Range("H32").Value = 1272.7551
$ws.Range("I32").Value = 1242.4045
$ws.Range("J32").Value = 1572.8889
$ws.Range("K32").Value = 1242.4045
$ws.Range("L32").Value = 1572.8889
$ws.Range("M32").Value = -955.4045000000001
$ws.Range("N32").Value = -2146.8889
# Row 61
$ws.Range("H61").Value = 1134.4375
$ws.Range("I61").Value = 1156.7333
$ws.Range("J61").Value = 800
$ws.Range("K61").Value = 1156.7333
$ws.Range("L61").Value = 800
$ws.Range("M61").Value = -944.7333000000001
$ws.Range("N61").Value = -1224
# Row 132
$ws.Range("H132").Value = 2354866.5
$ws.Range("I132").Value = 1425.5883
$ws.Range("J132").Value = 7355928
$ws.Range("K132").Value = 4276.7649
$ws.Range("L132").Value = 22067784
$ws.Range("M132").Value = -1746.7649
$ws.Range("N132").Value = -22072844
# Row 136
$ws.Range("H136").Value = 1134.4375
$ws.Range("I136").Value = 1156.7333
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 3470.199900000001
$ws.Range("L136").Value = 2400
$ws.Range("M136").Value = -920.1999000000005
$ws.Range("N136").Value = -7500

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 11822.1
$ws.Range("I33").Value = 7603
$ws.Range("J33").Value = 21666.666
$ws.Range("K33").Value = 7603
$ws.Range("L33").Value = 21666.666
$ws.Range("M33").Value = -7267
$ws.Range("N33").Value = -22338.666
# Row 82
$ws.Range("H82").Value = 14460.25
$ws.Range("I82").Value = 4852.6665
$ws.Range("J82").Value = 43283
$ws.Range("K82").Value = 4852.6665
$ws.Range("L82").Value = 43283
$ws.Range("M82").Value = -4469.6665
$ws.Range("N82").Value = -44049
# Row 85
$ws.Range("H85").Value = 14460.25
$ws.Range("I85").Value = 4852.6665
$ws.Range("J85").Value = 43283
$ws.Range("K85").Value = 4852.6665
$ws.Range("L85").Value = 43283
$ws.Range("M85").Value = -3526.6665
$ws.Range("N85").Value = -45935
# Row 97
$ws.Range("H97").Value = 2000
$ws.Range("I97").Value = 2000
$ws.Range("K97").Value = 2000
$ws.Range("M97").Value = -1009
# Row 102
$ws.Range("H102").Value = 17194.25
$ws.Range("I102").Value = 6990.3335
$ws.Range("J102").Value = 47806
$ws.Range("K102").Value = 6990.3335
$ws.Range("L102").Value = 47806
$ws.Range("M102").Value = -3745.3335
$ws.Range("N102").Value = -54296
# Row 134
$ws.Range("H134").Value = 2528594.2
$ws.Range("I134").Value = 792.625
$ws.Range("J134").Value = 9269399
$ws.Range("K134").Value = 2377.875
$ws.Range("L134").Value = 27808197
$ws.Range("M134").Value = 157.125
$ws.Range("N134").Value = -27813267

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 99
$ws.Range("H99").Value = 100003304
$ws.Range("I99").Value = 333336000
$ws.Range("J99").Value = 3571.4285
$ws.Range("K99").Value = 333336000
$ws.Range("L99").Value = 3571.4285
$ws.Range("M99").Value = -333334502
$ws.Range("N99").Value = -6567.4285
# Row 126
$ws.Range("H126").Value = 100003304
$ws.Range("I126").Value = 333336000
$ws.Range("J126").Value = 3571.4285
$ws.Range("K126").Value = 1000008000
$ws.Range("L126").Value = 10714.2855
$ws.Range("M126").Value = -1000005530
$ws.Range("N126").Value = -15654.2855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2730
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2064
$ws.Range("N67").ClearContents()
# Row 68
$ws.Range("H68").Value = 7136.8667
$ws.Range("I68").Value = 405.8889
$ws.Range("J68").Value = 17233.334
$ws.Range("K68").Value = 1217.6667
$ws.Range("L68").Value = 51700.00199999999
$ws.Range("M68").Value = -406.6667
$ws.Range("N68").Value = -53322.00199999999
# Row 69
$ws.Range("H69").Value = 9400.733
$ws.Range("J69").Value = 10036.5
$ws.Range("L69").Value = 30109.5
$ws.Range("N69").Value = -31731.5
# Row 70
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2685
$ws.Range("N70").ClearContents()
# Row 71
$ws.Range("H71").Value = 7136.8667
$ws.Range("I71").Value = 405.8889
$ws.Range("J71").Value = 17233.334
$ws.Range("K71").Value = 3653.0001
$ws.Range("L71").Value = 155100.006
$ws.Range("M71").Value = 402.9999000000003
$ws.Range("N71").Value = -163212.006
# Row 72
$ws.Range("H72").Value = 9400.733
$ws.Range("J72").Value = 10036.5
$ws.Range("L72").Value = 90328.5
$ws.Range("N72").Value = -98440.5
# Row 73
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1908
$ws.Range("N73").ClearContents()
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
# Row 75
$ws.Range("H75").Value = 5082
$ws.Range("I75").Value = 656.5
$ws.Range("J75").Value = 9507.5
$ws.Range("K75").Value = 1969.5
$ws.Range("L75").Value = 28522.5
$ws.Range("M75").Value = -971.5
$ws.Range("N75").Value = -30518.5
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
# Row 78
$ws.Range("H78").Value = 5082
$ws.Range("I78").Value = 656.5
$ws.Range("J78").Value = 9507.5
$ws.Range("K78").Value = 5908.5
$ws.Range("L78").Value = 85567.5
$ws.Range("M78").Value = -916.5
$ws.Range("N78").Value = -95551.5
# Row 107
$ws.Range("H107").Value = 37041484
$ws.Range("I107").Value = 202.3077
$ws.Range("J107").Value = 71436960
$ws.Range("K107").Value = 606.9231
$ws.Range("L107").Value = 214310880
$ws.Range("M107").Value = 1313.0769
$ws.Range("N107").Value = -214314720
# Row 122
$ws.Range("H122").Value = 7816581
$ws.Range("I122").Value = 29412058
$ws.Range("J122").Value = 5451.085
$ws.Range("K122").Value = 264708522
$ws.Range("L122").Value = 49059.765
$ws.Range("M122").Value = -264706072
$ws.Range("N122").Value = -53959.765

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 5571.6665
$ws.Range("I99").Value = 1540
$ws.Range("J99").Value = 17666.666
$ws.Range("K99").Value = 1540
$ws.Range("L99").Value = 17666.666
$ws.Range("M99").Value = 706
$ws.Range("N99").Value = -22158.666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2468
$ws.Range("I7").Value = 2752
$ws.Range("J7").Value = 1900
$ws.Range("K7").Value = 2752
$ws.Range("L7").Value = 1900
$ws.Range("M7").Value = -2640
$ws.Range("N7").Value = -2124
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
# Row 126
$ws.Range("H126").Value = 2468
$ws.Range("I126").Value = 2752
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 8256
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -5786
$ws.Range("N126").Value = -10640
# Row 132
$ws.Range("H132").Value = 5518.1294
$ws.Range("I132").Value = 1589.9062
$ws.Range("J132").Value = 11231.909
$ws.Range("K132").Value = 4769.7186
$ws.Range("L132").Value = 33695.727
$ws.Range("M132").Value = -2239.7186
$ws.Range("N132").Value = -38755.727
# Row 136
$ws.Range("H136").Value = 30890544
$ws.Range("I136").Value = 4764640.5
$ws.Range("J136").Value = 142858700
$ws.Range("K136").Value = 14293921.5
$ws.Range("L136").Value = 428576100
$ws.Range("M136").Value = -14291371.5
$ws.Range("N136").Value = -428581200

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
# Row 132
$ws.Range("H132").Value = 26549.762
$ws.Range("I132").Value = 31202.885
$ws.Range("J132").Value = 11744.363
$ws.Range("K132").Value = 93608.655
$ws.Range("L132").Value = 35233.089
$ws.Range("M132").Value = -91078.655
$ws.Range("N132").Value = -40293.089
